$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C ("denomination") shifting taille/prix_tvac/path right by one.
$ws.Columns("C:C").Insert()
$ws.Columns("C:C").ColumnWidth = 6

# Header
$ws.Range("C1").Value = "denomination"

# Denomination values per row (the size extracted out of the article name)
$den = @{
  2  = "1/2N"
  3  = "3/4N"
  4  = "4/4N"
  5  = "1/2R"
  6  = "3/4R"
  7  = "4/4R"
  8  = "1/2 NM"
  9  = "3/4 NM"
  10 = "4/4 NM"
  11 = "1/2 NF"
  12 = "3/4 NF"
  13 = "4/4 NF"
  14 = "1/2 RM"
  15 = "3/4 RM"
  16 = "4/4 RM"
  17 = "1/2 RF"
  18 = "3/4 RF"
  19 = "4/4 RF"
  20 = "6/4 RF"
  21 = "1/2 N"
  22 = "3/4 N"
  23 = "4/4 N"
  24 = "1/2 R"
  25 = "3/4 R"
  26 = "4/4 R"
  27 = "1/2 N"
  28 = "3/4 N"
  29 = "4/4 N"
  30 = "1/2 R"
  31 = "3/4 R"
  32 = "4/4 R"
}

foreach ($row in 2..32) {
  $ws.Range("C$row").Value = $den[$row]
}

$ws.Range("C4").Select()
